$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data for columns I (I0) and J (IF) per row
$ijData = @{
    2 = @(5, 6)
    3 = @(8, 8)
    4 = @(6, 7)
    5 = @(6, 7)
    6 = @(6, 7)
    7 = @(6, 7)
    8 = @(6, 7)
    9 = @(7, 8)
    10 = @(5, 6)
    11 = @(5, 6)
    12 = @(8, 9)
    13 = @(5, 6)
    14 = @(7, 7)
    15 = @(6, 6)
    16 = @(8, 8)
    17 = @(2, 5)
    18 = @(8, 8)
    19 = @(6, 8)
    20 = @(6, 7)
    21 = @(7, 8)
    22 = @(6, 7)
    23 = @(8, 9)
    24 = @(3, 4)
    25 = @(5, 7)
    26 = @(2, 5)
    27 = @(7, 7)
    28 = @(2, 4)
    29 = @(8, 8)
    30 = @(7, 8)
    31 = @(10, 10)
    32 = @(4, 5)
    33 = @(6, 6)
    34 = @(7, 7)
    35 = @(6, 6)
    36 = @(7, 7)
    37 = @(7, 8)
    38 = @(6, 7)
    39 = @(9, 9)
    40 = @(6, 7)
    41 = @(6, 7)
    42 = @(5, 7)
    43 = @(5, 6)
    44 = @(4, 4)
    45 = @(7, 7)
    46 = @(7, 8)
    47 = @(7, 7)
}

foreach ($row in $ijData.Keys) {
    $vals = $ijData[$row]
    $ws.Cells.Item([int]$row, 9).Value = $vals[0]
    $ws.Cells.Item([int]$row, 10).Value = $vals[1]
}
